$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("D228").Value
Write-Host "type:" $v.GetType()
Write-Host "value:" $v
$v2 = $ws.Cells.Item(228,4).Value
Write-Host "value2:" $v2
